# Weekly update: a new price observation is inserted above the existing
# row 136 ("Poroto granado" / Comercializadora del Agro de Limarí), pushing
# the old rows 136-139 down to 137-140 (dimension grows from R139 to R140).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 136, shifting rows 136:139 down to 137:140.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new observation.
$ws.Range("A136").Value = 2
$ws.Range("B136").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C136").Value = "Coquimbo"
$ws.Range("D136").Value = 45239
$ws.Range("E136").Value = 4
$ws.Range("F136").Value = 100112030
$ws.Range("G136").Value = "Poroto granado"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 400
$ws.Range("K136").Value = 28000
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = 29000
$ws.Range("N136").Value = "$/caja 15 kilos"
$ws.Range("O136").Value = "Provincia de Limarí"
$ws.Range("P136").Value = 1933
$ws.Range("Q136").Value = 15
$ws.Range("R136").Value = "Hortaliza"
